# Prep predict3dunet-231129-0 session with sample data & model.
#
# Row 75 previously only held the session id ("231129-0") and a short
# description of the experiment. This edit fleshes the row out into a
# full test-case entry (description / expected result / actual result)
# and fills the remaining tracking columns with the "TBD" placeholder
# that is used throughout the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: session id stays "231129-0" (unchanged, just re-asserted) ---
$ws.Range("A75").Value = "231129-0"

# --- Column B: tool name stays "predict3dunet" (unchanged) ---
$ws.Range("B75").Value = "predict3dunet"

# --- Column C: expanded description of the experiment ---
$ws.Range("C75").Value = "Attempt using Wolny's pre-trained model to segment his sample data and sample config files from his README.md."

# --- Columns D & E are new for this row: expected result / success criteria ---
$ws.Range("D75").Value = "The segmentation will work without error. The quality does not matter."
$ws.Range("E75").Value = "TBD: Success if no error, failure if error."

# --- Remaining tracking columns (F..BC, skipping the hidden AJ column
#     which stays blank) default to the "TBD" placeholder, same as the
#     other freshly-added rows further down the sheet. ---
$ws.Range("F75:AI75").Value = "TBD"
$ws.Range("AK75:BC75").Value = "TBD"

# Re-apply the shared "Calibri / vertically centered" look used across
# this row to the cells that now hold real text (A, C, D, E, G), so the
# new text matches the appearance of the rest of the sheet.
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4122)
$ws.Range("C75").PasteSpecial(-4122)
$ws.Range("D75").PasteSpecial(-4122)
$ws.Range("E75").PasteSpecial(-4122)
$ws.Range("G75").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Leave the selection on A75, matching where the cursor ended up after
# this edit.
$ws.Range("A75").Select()
